# Update "想去人数" (F column) counts on the 展览 and 全部类型 sheets
# to reflect newly output data.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5717
$ws1.Range("F3").Value = 88
$ws1.Range("F5").Value = 967
$ws1.Range("F7").Value = 2632
$ws1.Range("F9").Value = 190
$ws1.Range("F13").Value = 2468
$ws1.Range("F14").Value = 515

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5717
$ws4.Range("F3").Value = 88
$ws4.Range("F6").Value = 967
$ws4.Range("F8").Value = 2632
$ws4.Range("F10").Value = 190
$ws4.Range("F15").Value = 2468
$ws4.Range("F16").Value = 515
